# TC-51390: Verify_Backplane_Logic_For_Pro32xD_Pro32xBB
# Updates several "NA" (shared-string) cells to boolean FALSE on both
# sheets, auto-fits row 13 on "Add Devices" (drops the explicit row
# height), and swaps which sheet/cell is the active selection so that
# "Device_With_Order_Different" becomes the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Add Devices")
$ws2 = $wb.Worksheets.Item("Device_With_Order_Different")

# --- "Add Devices" sheet: NA -> FALSE on rows 10-13 -------------------
$ws1.Range("L10").Value = $false
$ws1.Range("M10").Value = $false
$ws1.Range("M11").Value = $false
$ws1.Range("L12").Value = $false
$ws1.Range("M12").Value = $false
$ws1.Range("M13").Value = $false

# Row 13 loses its explicit (wrapped-text) height -> back to auto height.
$ws1.Rows.Item(13).AutoFit()

# --- "Device_With_Order_Different" sheet: NA -> FALSE on row 10 -------
$ws2.Range("L10").Value = $false
$ws2.Range("M10").Value = $false

# --- Selection / active-tab bookkeeping --------------------------------
# "Add Devices" is no longer the selected/active tab; its lingering
# selection moves to L14.
$ws1.Range("L14").Select()

# "Device_With_Order_Different" becomes the selected/active tab with its
# selection at A10.
$ws2.Range("A10").Select()
$ws2.Activate()
